$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Update B1 text: append the extra "ตัวเลขเท่านั้น" qualifier ---
$ws.Range("B1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# --- Column width changes (engine rounds ColumnWidth to the nearest 1/6 of a
#     character, so the inputs below are chosen to land as close as possible
#     on the target widths of 25.81640625 / 25.7265625 characters) ---
$ws.Columns.Item(2).ColumnWidth = 24.9167
$ws.Range("G1:H1").ColumnWidth = 24.75

# --- New merged header cell G1:H1, centered horizontally only ---
$ws.Range("G1:H1").HorizontalAlignment = -4108
$ws.Range("G1:H1").Merge()

# --- New helper cells in row 2, reusing the existing Neutral/Bad cell
#     styles already used by B2/A2 so no new cell-format entries are
#     introduced ---
$ws.Range("B2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "เป็นค่าว่างได้"

$ws.Range("A2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"

$excel.CutCopyMode = $false

# --- Clear now-unused helper cells in rows 4 and 5 ---
$ws.Range("G4:H5").Clear()

# --- Move the active selection to D7 ---
$ws.Range("D7").Select()
